# IHW.xlsx: nieuwe settings voor samenstellen OpenAPI
#
# Moves the per-profile (SIM/UGM/BSM) yes/no configuration matrix from the
# "createjsonschema"/"createxmlschema" rows to the "createyaml" row, and
# gives "createjsonschema"/"createxmlschema" (and the two *variant rows)
# a simple static "no"/"NVT" value in column F instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IHW")

# --- Row 29: createjsonschema ---------------------------------------
# Column F gets a static "no" value; the per-profile matrix (H:K, M:P, R:U)
# is cleared.
$ws.Range("F29").Value = "no"
$ws.Range("H29:K29").ClearContents()
$ws.Range("M29:P29").ClearContents()
$ws.Range("R29:U29").ClearContents()

# --- Row 30: createjsonschemavariant ---------------------------------
# Default value changes from "ihw" to the new "NVT" marker.
$ws.Range("F30").Value = "NVT"

# --- Row 39: createxmlschema ------------------------------------------
# Same treatment as row 29.
$ws.Range("F39").Value = "no"
$ws.Range("H39:K39").ClearContents()
$ws.Range("M39:P39").ClearContents()
$ws.Range("R39:U39").ClearContents()

# --- Row 40: createxmlschemavariant -----------------------------------
# New static "NVT" value in column F (previously empty).
$ws.Range("F40").Value = "NVT"

# --- Row 41: createyaml -------------------------------------------------
# Column F static "no" is removed; instead the row now carries the
# per-profile yes/no matrix that used to live on row 29/39.
$ws.Range("F41").ClearContents()

$ws.Range("H41").Value = "no"
$ws.Range("I41").Value = "no"
$ws.Range("J41").Value = "no"
$ws.Range("K41").Value = "no"

$ws.Range("M41").Value = "no"
$ws.Range("N41").Value = "no"
$ws.Range("O41").Value = "no"
$ws.Range("P41").Value = "no"

$ws.Range("R41").Value = "no"
$ws.Range("S41").Value = "yes"
$ws.Range("T41").Value = "yes"
$ws.Range("U41").Value = "yes"

# --- View state: leave the selection on the newly edited cell ----------
$ws.Activate() | Out-Null
$ws.Range("T41").Select() | Out-Null
